# "SQL Join parents E-Mail" sprint-2 burndown update
# Day column "D" (sprint day 2) gets filled in with effort-remaining figures
# for the first five backlog items, mirroring the pattern already used for
# column C (day 1). Rows 13-37 have no stories yet, so column D stays blank
# there, same as column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current Iteration")
$wsChart = $wb.Worksheets.Item("Burndown Chart")

# --- Column D (day 2) for rows 8-12 -----------------------------------
# Copy the "has data" cell style from column C onto column D for each row
# before writing the new formulas/values, so D matches the look of a
# filled-in day (fill 6) instead of the blank/unused day style (fill 7).
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D8").Formula = "=B8/(7+2)*6"
$ws.Range("D9").Value = 0
$ws.Range("D10").Formula = "=B10/(6.5+1+1)*5"
$ws.Range("D11").Formula = "=B11/(7+2+0)*6"
$ws.Range("D12").Formula = "=C12/(10+2+2)*5"

# Recalculate so Task Balance (M), Percentage Completed (N), and the day-2
# "Effort Remaining" row (38) all pick up the new column D figures.
$excel.CalculateFullRebuild()

# --- Restore navigation / selection state ------------------------------
$ws.Activate()
[void]$ws.Range("Q8").Select()

$wsChart.Activate()
[void]$wsChart.Range("C48").Select()

$ws.Activate()
